# Add chapter-2 quiz questions (German statement, English translation, True/False answer)
# to the rows immediately following the existing question bank (rows 1-26 already hold
# chapter-1 questions). New rows are appended starting at row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questions = @(
    @("Durch den Multiplikatoreffekt werden wirtschaftliche Schocks verstärkt.", "Economic shocks are amplified by the multiplier effect", $true),
    @("Durch den Multiplikatoreffekt werden wirtschaftliche Schocks abgeschwächt.", " Economic shocks are mitigated by the multiplier effect", $false),
    @("Importe sind Teil des Bruttoinlandsproduktes. ", " Imports are part of the gross domestic product", $false),
    @("Exporte sind Teil des Bruttoinlandsproduktes. ", " Exports are part of the gross domestic product", $true),
    @("In unserem Gütermarktmodell werden Lagerinvestitionen nicht berücksichtigt.", " In our goods market model, inventory investments are not considered", $true),
    @("In unserem Gütermarktmodell werden Lagerinvestitionen berücksichtigt.", " In our goods market model, inventory investments are considered", $false),
    @("In unserem Gütermarktmodell ignorieren wir Preisveränderungen.", " In our goods market model, we ignore price changes", $true),
    @("Der Konsum hängt in erster Linie vom verfügbaren Einkommen ab.", " Consumption primarily depends on disposable income", $true),
    @("Das verfügbare Einkommen ist als die Differenz zwischen Steuern und Einkommen definiert.", " Disposable income is defined as the difference between taxes and income", $false),
    @("Das verfügbare Einkommen ist als die Differenz zwischen Einkommen und Steuern definiert.", " Disposable income is defined as the difference between income and taxes", $true),
    @("Die marginale Konsumquote kann größer als 1 sein.", " The marginal propensity to consume can be greater than 1", $false),
    @("Die marginale Konsumquote kann nicht größer als 1 sein.", " The marginal propensity to consume cannot be greater than 1", $true),
    @("Entscheidungen über die Höhe der Staatsausgaben und der Steuern bezeichnet man als Fiskalpolitik.", " Decisions about the level of government spending and taxes are called fiscal policy", $true),
    @("Entscheidungen über die Höhe der Staatsausgaben und der Steuern bezeichnet man als Steuerpolitik.", " Decisions about the level of government spending and taxes are called tax policy", $false),
    @("Eine Variable, die von anderen Variablen abhängig ist, nennt man endogen.", " A variable that depends on other variables is called endogenous", $true),
    @("Eine Variable, die nicht von den anderen Variablen eines Modells abhängt, nennt man exogen.", " A variable that does not depend on other variables in a model is called exogenous", $true),
    @("Eine Variable, die von anderen Variablen abhängig ist, nennt man exogen.", " A variable that depends on other variables is called exogenous", $false),
    @("Eine Variable, die nicht von den anderen Variablen eines Modells abhängt, nennt man endogen.", " A variable that does not depend on other variables in a model is called endogenous", $false),
    @("Im Gütermarktgleichgewicht ist die Produktion von Gütern gleich der Nachfrage von Gütern.", " In goods market equilibrium, the production of goods equals the demand for goods", $true),
    @("Im Gütermarktgleichgewicht ist die Produktion von Gütern gleich dem Einkommen.", " In goods market equilibrium, the production of goods equals income", $true),
    @("Im Gütermarktgleichgewicht ist die Produktion von Gütern gleich dem Staatsausgaben.", " In goods market equilibrium, the production of goods equals government spending", $false),
    @("Im Gütermarktgleichgewicht sind die Investionen gleich dem öffentlichen Sparen.", " In goods market equilibrium, investments equal public savings", $false),
    @("Im Gütermarktgleichgewicht sind die Investionen gleich der gesamtwirtschaftlichen Ersparnis.", " In goods market equilibrium, investments equal aggregate savings", $true),
    @("Einkommensabhängige Steuern verringern den Multiplikator. ", " Income-dependent taxes reduce the multiplier", $true),
    @("Einkommensabhängige Steuern erhöhen den Multiplikator. ", " Income-dependent taxes increase the multiplier", $false),
    @("Eine Erhöhung der Staatsausgaben führt im Gütermarktmodell zu einer Verschiebung der Nachfragekurve nach oben.", " An increase in government spending leads to an upward shift in the demand curve in the goods market model", $true),
    @("Eine Senkung der Staatsausgaben führt im Gütermarktmodell zu einer Verschiebung der Nachfragekurve nach unten.", " A decrease in government spending leads to a downward shift in the demand curve in the goods market model", $true),
    @("Eine Erhöhung der Steuern führt im Gütermarktmodell zu einer Verschiebung der Nachfragekurve nach oben.", " An increase in taxes leads to an upward shift in the demand curve in the goods market model", $false),
    @("Eine Senkung der Steuern führt im Gütermarktmodell zu einer Verschiebung der Nachfragekurve nach unten.", " A decrease in taxes leads to a downward shift in the demand curve in the goods market model", $false),
    @("Eine Erhöhung der Staatsausgaben führt im Gütermarktmodell zu einer Verschiebung der Nachfragekurve nach unten.", " An increase in government spending leads to a downward shift in the demand curve in the goods market model", $true),
    @("Eine Senkung der Staatsausgaben führt im Gütermarktmodell zu einer Verschiebung der Nachfragekurve nach oben.", " A decrease in government spending leads to an upward shift in the demand curve in the goods market model", $true),
    @("Eine Erhöhung der Steuern führt im Gütermarktmodell zu einer Verschiebung der Nachfragekurve nach unten.", " An increase in taxes leads to a downward shift in the demand curve in the goods market model", $false),
    @("Eine Senkung der Steuern führt im Gütermarktmodell zu einer Verschiebung der Nachfragekurve nach oben.", " A decrease in taxes leads to an upward shift in the demand curve in the goods market model", $false),
    @("Steigt im Gütermarktmodell die marginale Konsumquote, wird die Nachfragekurve steiler. ", " If the marginal propensity to consume rises in the goods market model, the demand curve becomes steeper", $true),
    @("Steigt im Gütermarktmodell die marginale Konsumquote, wird die Nachfragekurve flacher. ", " If the marginal propensity to consume rises in the goods market model, the demand curve becomes flatter.", $false)
)

$startRow = 27
for ($i = 0; $i -lt $questions.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $questions[$i][0]
    $ws.Cells.Item($row, 2).Value = $questions[$i][1]
    $ws.Cells.Item($row, 3).Value = $questions[$i][2]
}

# Cosmetic touch-ups matching the author's final view state: column A keeps its
# auto-fit-to-content width, column B is narrowed to a fixed width (no longer
# auto-fit), and the active selection/scroll position moves down to the newly
# added question block.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).ColumnWidth = 34.29

$ws.Range("A36").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 26
$win.ScrollColumn = 1
